$wb = $excel.ActiveWorkbook

# --- 1. Rename "Data" -> "Storm Event 1" ---
$ws1 = $wb.Worksheets.Item("Data")
$ws1.Name = "Storm Event 1"

# --- 2. Update the instructional text in column E (reworded + extra sentence) ---
$ws1.Range("E2").Value = "Delete this row too. It's only here for your preference what datatype that the app expects for each column. Please DO NOT delete/rename the datetime column"

# --- 3. Widen column E to fit the longer instructional text ---
$ws1.Columns.Item(5).ColumnWidth = 145.45

# --- 4. Add a second data sheet "Storm Event 2" after "Storm Event 1" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Storm Event 2"

$ws2.Range("A1").Value = "datetime"
$ws2.Range("B1").Value = "Sensor A"
$ws2.Range("C1").Value = "Sensor B"
$ws2.Range("D1").Value = "Sensor C"
$ws2.Range("A2").Value = "2021-12-14 05:14:00"
$ws2.Range("B2").Value = 50.17
$ws2.Range("C2").Value = 24.88
$ws2.Range("D2").Value = 23.75
$ws2.Range("E2").Value = "Delete this row too. It's only here for your preference what datatype that the app expects for each column. Please DO NOT delete/rename the datetime column"

# match column layout/widths of "Storm Event 1"
$ws2.Columns.Item(1).ColumnWidth = $ws1.Columns.Item(1).ColumnWidth
$ws2.Columns.Item(2).ColumnWidth = $ws1.Columns.Item(2).ColumnWidth
$ws2.Columns.Item(3).ColumnWidth = $ws1.Columns.Item(3).ColumnWidth
$ws2.Columns.Item(4).ColumnWidth = $ws1.Columns.Item(4).ColumnWidth
$ws2.Columns.Item(5).ColumnWidth = 145.45

# match header / row styling of "Storm Event 1"
$ws2.Range("A1:D1").Style = $ws1.Range("A1:D1").Style
$ws2.Range("A2:E2").Style = $ws1.Range("A2:E2").Style

# --- 5. Update sheet view / selection state on each sheet ---
$wsInstr = $wb.Worksheets.Item("Instructions")
$wsInstr.Activate()
$wsInstr.Range("A2").Select()

$ws1.Activate()
$ws1.Range("B3").Select()
$ws1.Cells.Select()

$ws2.Activate()
$ws2.Range("D30").Select()
